$d = $word.ActiveDocument

# --- Heading text: insert a space between the word and the milestone/activity number ---
# "MilestoneN:" -> "Milestone N:" (also handles "Milestone5:Project" -> "Milestone 5:Project")
$d.Content.Find.Execute("Milestone1:", $true, $false, $false, $false, $false, $true, 1, $false, "Milestone 1:", 2)
$d.Content.Find.Execute("Milestone2:", $true, $false, $false, $false, $false, $true, 1, $false, "Milestone 2:", 2)
$d.Content.Find.Execute("Milestone5:", $true, $false, $false, $false, $false, $true, 1, $false, "Milestone 5:", 2)
$d.Content.Find.Execute("Milestone6:", $true, $false, $false, $false, $false, $true, 1, $false, "Milestone 6:", 2)

# "Milestone 5:Project" -> "Milestone 5: Project" (space before "Project")
$d.Content.Find.Execute("Milestone 5:Project", $true, $false, $false, $false, $false, $true, 1, $false, "Milestone 5: Project", 2)

# "ActivityN" -> "Activity N" (covers both "ActivityN:" and "ActivityN" forms, colon may be a
# separate run so match on the word+digit only)
$d.Content.Find.Execute("Activity1", $true, $false, $false, $false, $false, $true, 1, $false, "Activity 1", 2)
$d.Content.Find.Execute("Activity2", $true, $false, $false, $false, $false, $true, 1, $false, "Activity 2", 2)
$d.Content.Find.Execute("Activity3", $true, $false, $false, $false, $false, $true, 1, $false, "Activity 3", 2)

# --- Body text fixes ---

# "... identifying different water bottle images .It encompasses ..."
#   -> "... identifying different water bottle images. It encompasses ..."
$d.Content.Find.Execute("different water bottle images .It", $true, $false, $false, $false, $false, $true, 1, $false, "different water bottle images. It", 2)

# "... effectiveness of subsequent analyses ." (end of paragraph) -> "... subsequent analyses."
$d.Content.Find.Execute("effectiveness of subsequent analyses .", $true, $false, $false, $false, $false, $true, 1, $false, "effectiveness of subsequent analyses.", 2)

# First CNN model list (uses British "optimiser" - unique to this paragraph):
# "CNN (batch normalisation and Adam), CNN (adam optimizer 2 fully connected layers), CNN (with
#  optimiser adam ), CNN (with SDG) )" ->
# "CNN (batch normalization and Adam), CNN (adam optimizer 2 fully connected layers), CNN (with
#  optimizer adam), CNN (with SDG))"
$d.Content.Find.Execute("CNN (with optimiser adam ), CNN (with SDG) )", $true, $false, $false, $false, $false, $true, 1, $false, "CNN (with optimizer adam), CNN (with SDG))", 2)
$d.Content.Find.Execute("CNN (batch normalisation and Adam), CNN (adam optimizer 2 fully connected layers), CNN (with optimizer adam)", $true, $false, $false, $false, $false, $true, 1, $false, "CNN (batch normalization and Adam), CNN (adam optimizer 2 fully connected layers), CNN (with optimizer adam)", 2)

# Second CNN model list (American "optimizer" already, ends with ". It considers"):
# "CNN (batch normalisation and Adam), CNN (adam optimizer 2 fully connected layers), CNN (with
#  optimizer adam ), CNN (with SDG) . It considers" ->
# "CNN (batch normalization and adam), CNN (adam optimizer 2 fully connected layers), CNN (with
#  optimizer adam), CNN (with SDG). It considers"
$d.Content.Find.Execute("CNN (with optimizer adam ), CNN (with SDG) . It considers", $true, $false, $false, $false, $false, $true, 1, $false, "CNN (with optimizer adam), CNN (with SDG). It considers", 2)
$d.Content.Find.Execute("CNN (batch normalisation and Adam), CNN (adam optimizer 2 fully connected layers), CNN (with optimizer adam)", $true, $false, $false, $false, $false, $true, 1, $false, "CNN (batch normalization and adam), CNN (adam optimizer 2 fully connected layers), CNN (with optimizer adam)", 2)

# "modeling.The subsequent" -> "modeling. The subsequent"
$d.Content.Find.Execute("modeling.The subsequent", $true, $false, $false, $false, $false, $true, 1, $false, "modeling. The subsequent", 2)

# "... file submission in Github, Kindly click ..." -> "... file submission in GitHub, kindly click ..."
$d.Content.Find.Execute("Github,", $true, $false, $false, $false, $false, $true, 1, $false, "GitHub,", 2)
$d.Content.Find.Execute("Kindly", $true, $false, $false, $false, $false, $true, 1, $false, "kindly", 2)
